# Duplicate the "status" sheet, placing the copy immediately before it,
# then rename the copy to "action" and make it the active sheet/tab,
# matching the "good separation utility and work order" edit:
#   - a new "action" worksheet is inserted right before "status"
#   - it carries the same table (name/description lookup rows) as "status"
#   - it becomes the active sheet (selected cell C28)
#   - the previously active "User" sheet is left on cell B27, no longer selected

$wb = $excel.ActiveWorkbook

$statusSheet = $wb.Worksheets.Item("status")
$statusSheet.Copy($statusSheet)
$actionSheet = $wb.ActiveSheet
$actionSheet.Name = "action"

$userSheet = $wb.Worksheets.Item("User")
$userSheet.Activate()
$userSheet.Range("B27").Select() | Out-Null

$actionSheet.Activate()
$actionSheet.Range("C28").Select() | Out-Null
